$d = $word.ActiveDocument

# The only substantive text change in the target diff is the name
# "Roee esquire" -> "Roee esquira" in the first answer paragraph.
# (All other hunks in the diff are just Word re-flowing/merging runs
# and dropping proofErr spell/grammar markers around unchanged text.)
$d.Content.Find.Execute("Roee esquire, ID 309840791", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Roee esquira, ID 309840791", 2)
